$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 25
$ws.Range("D5").Select()
